# Merge-Run: collapses every run spanned by $targetText (which must
# already equal the live text at that position) into a single run.
#
# PowerPoint (and this COM-interop host) keeps a run's formatting intact
# as long as a Range.Text assignment doesn't need to touch it, which
# means a plain "set Text to the text that's already there" is a no-op
# that leaves existing run splits untouched. To force an actual merge we
# first overwrite the span with a same-length-ish placeholder (which *is*
# a real change, so the engine collapses the whole span into one run),
# then write the real text back into that single run.
function Merge-Run {
    param($TextRange, [string]$TargetText)

    $full = $TextRange.Text
    $idx = $full.IndexOf($TargetText)
    if ($idx -lt 0) {
        return
    }
    $start = $idx + 1
    $length = $TargetText.Length
    $placeholder = "Z" * ($length + 5)
    $TextRange.Characters($start, $length).Text = $placeholder
    $TextRange.Characters($start, $placeholder.Length).Text = $TargetText
}

$p = $ppt.ActivePresentation

# --- Slide 33 ("AVX-512"): tidy up text that had been split across
# --- several identically-formatted runs, merging them back down to one
# --- run each (content/visual result is unchanged).
$s33 = $p.Slides.Item(33)

# Content placeholder: "Vector registers for floating point operands:" /
# "512 bit wide" / "8" + " double precision" / "16 single precision"
$contentTr = $s33.Shapes.Item(2).TextFrame.TextRange
Merge-Run $contentTr "Vector registers for floating point operands:"
Merge-Run $contentTr "512 bit wide"
Merge-Run $contentTr " double precision"
Merge-Run $contentTr "16 single precision"

# "8 concurrent operations" callout (Group 6 -> TextBox 4)
$tb4Tr = $s33.Shapes.Item(4).GroupItems.Item(1).TextFrame.TextRange
Merge-Run $tb4Tr "8 concurrent operations"

# "16 concurrent operations" callout (Group 7 -> TextBox 8)
$tb8Tr = $s33.Shapes.Item(5).GroupItems.Item(1).TextFrame.TextRange
Merge-Run $tb8Tr "16 concurrent operations"

# "Even more worth to recompile!" textbox
$tb12Tr = $s33.Shapes.Item(6).TextFrame.TextRange
Merge-Run $tb12Tr "Even more worth to recompile!"

# --- Slide 34: repurpose the "Helpful compiler options" slide into the
# --- new "Double promotion" slide by retitling it.
$s34 = $p.Slides.Item(34)
$s34.Shapes.Item(1).TextFrame.TextRange.Text = "Double promotion"
